# Generate Report for Handoff
#
# The localization handback/status report is being regenerated: the status
# moves from "Handed back: in sync with en-US" to "Ready for handoff", and
# the associated timestamps advance a little (the report was re-run).

$wb = $excel.ActiveWorkbook

$newStatus       = "Ready for handoff"
$newHoDate       = "2016-08-26 22:58:27"
$newHandoffDate  = "2016-08-26 22:58:22"

# --- Overview sheet ---------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("G2").Value = $newHoDate

# The status text got a lot shorter, so the "zh-cn"/"de-de" columns on the
# Overview sheet shrink to fit the new content.
$wsOverview.Columns("E").ColumnWidth = 16.35
$wsOverview.Columns("F").ColumnWidth = 16.35

# --- zh-cn sheet --------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("H2").Value = $newHandoffDate

# The "Status" column shrinks to fit the shorter text too.
$wsZhCn.Columns("C").ColumnWidth = 16.35

# --- de-de sheet --------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("H2").Value = $newHoDate

$wsDeDe.Columns("C").ColumnWidth = 16.35
